$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" property value (row 8, column B)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-04-30T15:08:51+00:00"

# ---------------------------------------------------------------------------
# 2. Mapping Table 0 sheet: add the new ConceptMap mapping rows
# ---------------------------------------------------------------------------
$map = $wb.Worksheets.Item("Mapping Table 0")

# Fill in the "Target" column (D) for the rows that already exist (rows 3-12)
$map.Range("D3").Value  = "DocumentReference.id"
$map.Range("D4").Value  = "DocumentReference.identifier"
$map.Range("D5").Value  = "DocumentReference.content.attachment.hash"
$map.Range("D6").Value  = "DocumentReference.content.attachment.size"
$map.Range("D7").Value  = "DocumentReference.content.attachment.language"
$map.Range("D8").Value  = "DocumentReference.authenticator"
$map.Range("D9").Value  = "DocumentReference.context.period.start"
$map.Range("D10").Value = "DocumentReference.context.period.end"
$map.Range("D11").Value = "DocumentReference.subject.fr-core-patient"
$map.Range("D12").Value = "DocumentReference.subject.fr-core-patient"

# Re-point the "Source" column (A) so every existing row shifts down to make
# room for the two brand-new concepts (entryUUID, logicalId) that now lead
# the table.
$map.Range("A3").Value  = "DocumentEntry.entryUUID"
$map.Range("A4").Value  = "DocumentEntry.logicalId"
$map.Range("A5").Value  = "DocumentEntry.hash"
$map.Range("A6").Value  = "DocumentEntry.size"
$map.Range("A7").Value  = "DocumentEntry.languageCode"
$map.Range("A8").Value  = "DocumentEntry.legalAuthenticator"
$map.Range("A9").Value  = "DocumentEntry.serviceStartTime"
$map.Range("A10").Value = "DocumentEntry.serviceEndTime"
$map.Range("A11").Value = "DocumentEntry.sourcePatientID"
$map.Range("A12").Value = "DocumentEntry.sourcePatientInfo"
$map.Range("A13").Value = "DocumentEntry.URI"
$map.Range("A14").Value = "DocumentEntry.title"
$map.Range("A15").Value = "DocumentEntry.comments"
$map.Range("A16").Value = "DocumentEntry.patientID"
$map.Range("A17").Value = "DocumentEntry.uniqueId"
$map.Range("A18").Value = "DocumentEntry.class"
$map.Range("A19").Value = "DocumentEntry.confidentiality"

# Append the rows that used to fall off the end of the table, plus the brand
# new rows 20-29, copying the formatting from the last pre-existing data row
# (19) so every new cell keeps style "s=2" like its neighbours.
$fmtSrc = $map.Range("A19:E19")

$newRows = @(
    "DocumentEntry.eventCodeList",
    "DocumentEntry.format",
    "DocumentEntry.healthcareFacilityTypeCode",
    "DocumentEntry.practiceSetting",
    "DocumentEntry.type",
    "DocumentEntry.documentAvailability",
    "DocumentEntry.homeCommunityId",
    "DocumentEntry.creationTime",
    "DocumentEntry.referenceIdList",
    "DocumentEntry.referenceIdList"
)

$rowIndex = 20
foreach ($source in $newRows) {
    $fmtSrc.Copy()
    $map.Range("A$rowIndex`:E$rowIndex").PasteSpecial(-4122)

    $map.Range("A$rowIndex").Value = $source
    $map.Range("C$rowIndex").Value = "equivalent"

    $rowIndex = $rowIndex + 1
}

$excel.CutCopyMode = 0
